$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "62.795.12"
$ws.Cells.Item(2, 5).Value = "  +3.09%  "
$ws.Cells.Item(3, 4).Value = "3.445.12"
$ws.Cells.Item(3, 5).Value = "  +2.04%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "582.16"
$ws.Cells.Item(5, 5).Value = "  +2.33%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "146.97"
$ws.Cells.Item(6, 5).Value = "  +5.13%  "
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.477"
$ws.Cells.Item(8, 5).Value = "  +1.05%  "
$ws.Cells.Item(9, 5).Value = "  +0.68%  "
$ws.Cells.Item(10, 5).Value = "  +2.79%  "
$ws.Cells.Item(11, 5).Value = "  +1.87%  "
$ws.Cells.Item(12, 4).Value = "4.035.94"
$ws.Cells.Item(12, 5).Value = "  +2.10%  "
$ws.Cells.Item(13, 5).Value = "  +5.21%  "
$ws.Cells.Item(14, 5).Value = "  -0.69%  "
$ws.Cells.Item(15, 4).Value = "3.442.95"
$ws.Cells.Item(15, 5).Value = "  +1.80%  "
$ws.Cells.Item(16, 5).Value = "  +2.45%  "
$ws.Cells.Item(17, 4).Value = "62.778.95"
$ws.Cells.Item(17, 5).Value = "  +2.85%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.23"
$ws.Cells.Item(18, 5).Value = "  +2.59%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.28"
$ws.Cells.Item(19, 5).Value = "  +5.68%  "
$ws.Cells.Item(20, 5).Value = "  +5.16%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "395.78"
$ws.Cells.Item(21, 5).Value = "  +4.03%  "
$ws.Cells.Item(22, 2).Value = "Litecoin"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "75.61"
$ws.Cells.Item(22, 5).Value = "  +0.31%  "
$ws.Cells.Item(23, 2).Value = "Polygon"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.562"
$ws.Cells.Item(23, 5).Value = "  +2.56%  "
$ws.Cells.Item(24, 5).Value = "  -0.01%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.0000118"
$ws.Cells.Item(25, 5).Value = "  +4.68%  "
$ws.Cells.Item(26, 4).Value = "3.584.88"
$ws.Cells.Item(26, 5).Value = "  +1.92%  "
$ws.Cells.Item(27, 5).Value = "  -0.99%  "
$ws.Cells.Item(28, 5).Value = "  +7.11%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.999"
$ws.Cells.Item(29, 5).Value = "  -0.09%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.18"
$ws.Cells.Item(30, 5).Value = "  +3.15%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.14"
$ws.Cells.Item(31, 5).Value = "  +1.07%  "
$ws.Cells.Item(32, 5).Value = "  +4.80%  "
$ws.Cells.Item(33, 5).Value = "  +0.00%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "23.82"
$ws.Cells.Item(34, 5).Value = "  +2.98%  "
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.62"
$ws.Cells.Item(35, 5).Value = "  +11.66%  "
$ws.Cells.Item(36, 2).Value = "NEARProtocol"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.31"
$ws.Cells.Item(36, 5).Value = "  +7.59%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "7.06"
$ws.Cells.Item(37, 5).Value = "  +2.41%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "168.59"
$ws.Cells.Item(38, 5).Value = "  +1.43%  "
$ws.Cells.Item(39, 4).Value = "3.479.09"
$ws.Cells.Item(39, 5).Value = "  +1.99%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "30.15"
$ws.Cells.Item(40, 5).Value = "  +16.46%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0769"
$ws.Cells.Item(41, 5).Value = "  +0.99%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.791"
$ws.Cells.Item(42, 5).Value = "  +1.60%  "
$ws.Cells.Item(43, 2).Value = "OKB"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "42.97"
$ws.Cells.Item(43, 5).Value = "  +1.55%  "
$ws.Cells.Item(44, 2).Value = "Filecoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "4.48"
$ws.Cells.Item(44, 5).Value = "  +3.25%  "
$ws.Cells.Item(45, 2).Value = "Stacks"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.72"
$ws.Cells.Item(45, 5).Value = "  +5.40%  "
$ws.Cells.Item(46, 2).Value = "ONDO"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.19"
$ws.Cells.Item(46, 5).Value = "  +7.64%  "
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(47, 4).Value = "2.522.79"
$ws.Cells.Item(47, 5).Value = "  +3.99%  "
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "23.50"
$ws.Cells.Item(48, 5).Value = "  +2.49%  "
$ws.Cells.Item(49, 2).Value = "Cosmos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "6.72"
$ws.Cells.Item(49, 5).Value = "  +1.66%  "
$ws.Cells.Item(50, 2).Value = "dogwifhat"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.19"
$ws.Cells.Item(50, 5).Value = "  +5.21%  "
$ws.Cells.Item(51, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Cells.Item(51, 5).Value = "  -0.08%  "

Write-Output "Applied all cell updates"